# =====================================================================
# Scheduled market-data refresh for the Leve profit tracker workbook.
#
# Each worksheet tab (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) lists Leve
# (Levequest) turn-in items alongside live Market Board pricing pulled
# from the Universalis API:
#   H = currentAveragePrice     I = currentAveragePriceNQ
#   J = currentAveragePriceHQ   K = LevePriceNQ
#   L = LevePriceHQ             M = LeveProfitNQ   N = LeveProfitHQ
#
# This run refreshes those columns with newly-fetched prices for the
# rows whose market data changed since the last sync. Cells that did
# not previously exist (no quantity priced yet) are created, and cells
# that no longer have a meaningful value are cleared to match the
# refreshed source feed.
# =====================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# ALC tab
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")
# Row 8
$ws.Range("H8").Value = 1365
$ws.Range("I8").Value = 1099.2307
$ws.Range("K8").Value = 3297.6921
$ws.Range("M8").Value = -3158.6921
# Row 17
$ws.Range("H17").Value = 9443.691999999999
$ws.Range("J17").Value = 10105.667
$ws.Range("L17").Value = 30317.001
$ws.Range("N17").Value = -30653.001
# Row 69
$ws.Range("H69").Value = 7669.4707
$ws.Range("J69").Value = 7779.3125
$ws.Range("L69").Value = 23337.9375
$ws.Range("N69").Value = -25085.9375
# Row 72
$ws.Range("H72").Value = 7669.4707
$ws.Range("J72").Value = 7779.3125
$ws.Range("L72").Value = 70013.8125
$ws.Range("N72").Value = -78749.8125
# Row 76
$ws.Range("H76").Value = 6315.7896
$ws.Range("I76").Value = 5054.222
$ws.Range("J76").Value = 7451.2
$ws.Range("K76").Value = 5054.222
$ws.Range("L76").Value = 7451.2
$ws.Range("M76").Value = -4739.222
$ws.Range("N76").Value = -8081.2
# Row 79
$ws.Range("H79").Value = 6315.7896
$ws.Range("I79").Value = 5054.222
$ws.Range("J79").Value = 7451.2
$ws.Range("K79").Value = 5054.222
$ws.Range("L79").Value = 7451.2
$ws.Range("M79").Value = -3962.222
$ws.Range("N79").Value = -9635.200000000001
# Row 138
$ws.Range("H138").Value = 2871.8433
$ws.Range("I138").Value = 1959.375
$ws.Range("J138").Value = 3089.7463
$ws.Range("K138").Value = 5878.125
$ws.Range("L138").Value = 9269.2389
$ws.Range("M138").Value = -738.125
$ws.Range("N138").Value = -19549.2389

# ---------------------------------------------------------------
# ARM tab
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")
# Row 19
$ws.Range("H19").Value = 3501.8333
$ws.Range("I19").Value = 2750.25
$ws.Range("J19").Value = 5005
$ws.Range("K19").Value = 2750.25
$ws.Range("L19").Value = 5005
$ws.Range("M19").Value = -2521.25
$ws.Range("N19").Value = -5463
# Row 32
$ws.Range("H32").Value = 1984.375
$ws.Range("I32").Value = 1503.2467
$ws.Range("K32").Value = 1503.2467
$ws.Range("M32").Value = -1216.2467
# Row 61
$ws.Range("H61").Value = 5213.567
$ws.Range("I61").Value = 4195.4546
$ws.Range("K61").Value = 4195.4546
$ws.Range("M61").Value = -3983.4546
# Row 110
$ws.Range("H110").Value = 2192.3428
$ws.Range("I110").Value = 1627.2
$ws.Range("J110").Value = 5583.2
$ws.Range("K110").Value = 1627.2
$ws.Range("L110").Value = 5583.2
$ws.Range("M110").Value = 417.8
$ws.Range("N110").Value = -9673.200000000001
# Row 136
$ws.Range("H136").Value = 5213.567
$ws.Range("I136").Value = 4195.4546
$ws.Range("K136").Value = 12586.3638
$ws.Range("M136").Value = -10036.3638

# ---------------------------------------------------------------
# CRP tab
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 34033.89
$ws.Range("I31").Value = 4831.857
$ws.Range("K31").Value = 4831.857
$ws.Range("M31").Value = -4536.857
# Row 34
$ws.Range("H34").Value = 34033.89
$ws.Range("I34").Value = 4831.857
$ws.Range("K34").Value = 4831.857
$ws.Range("M34").Value = -4629.857
# Row 58
$ws.Range("H58").Value = 7233.2666
$ws.Range("I58").Value = 2213.5715
$ws.Range("J58").Value = 11625.5
$ws.Range("K58").Value = 2213.5715
$ws.Range("L58").Value = 11625.5
$ws.Range("M58").Value = -2010.5715
$ws.Range("N58").Value = -12031.5
# Row 93
$ws.Range("H93").Value = 24599.25
$ws.Range("I93").Value = 24599.25
$ws.Range("K93").Value = 24599.25
$ws.Range("M93").Value = -22727.25
# Row 97
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()
# Row 100
$ws.Range("H100").Value = 34000
$ws.Range("J100").Value = 34000
$ws.Range("L100").Value = 34000
$ws.Range("N100").Value = -36164
# Row 109
$ws.Range("H109").Value = 71730.25
$ws.Range("J109").Value = 71730.25
$ws.Range("L109").Value = 71730.25
$ws.Range("N109").Value = -73810.25
# Row 136
$ws.Range("H136").Value = 7233.2666
$ws.Range("I136").Value = 2213.5715
$ws.Range("J136").Value = 11625.5
$ws.Range("K136").Value = 6640.7145
$ws.Range("L136").Value = 34876.5
$ws.Range("M136").Value = -4090.7145
$ws.Range("N136").Value = -39976.5

# ---------------------------------------------------------------
# CUL tab
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")
# Row 22
$ws.Range("H22").Value = 2398.2
$ws.Range("J22").Value = 4103
$ws.Range("L22").Value = 12309
$ws.Range("N22").Value = -12647
# Row 24
$ws.Range("H24").Value = 400.5
$ws.Range("J24").Value = 400
$ws.Range("L24").Value = 1200
$ws.Range("N24").Value = -1660
# Row 25
$ws.Range("H25").Value = 182.125
$ws.Range("I25").Value = 151
$ws.Range("J25").Value = 400
$ws.Range("K25").Value = 453
$ws.Range("L25").Value = 1200
$ws.Range("M25").Value = -284
$ws.Range("N25").Value = -1538
# Row 27
$ws.Range("H27").Value = 2398.2
$ws.Range("J27").Value = 4103
$ws.Range("L27").Value = 12309
$ws.Range("N27").Value = -12513
# Row 30
$ws.Range("H30").Value = 182.125
$ws.Range("I30").Value = 151
$ws.Range("J30").Value = 400
$ws.Range("K30").Value = 453
$ws.Range("L30").Value = 1200
$ws.Range("M30").Value = -351
$ws.Range("N30").Value = -1404
# Row 34
$ws.Range("H34").Value = 2773.75
$ws.Range("I34").Value = 1784.125
$ws.Range("J34").Value = 3433.5
$ws.Range("K34").Value = 5352.375
$ws.Range("L34").Value = 10300.5
$ws.Range("M34").Value = -5268.375
$ws.Range("N34").Value = -10468.5
# Row 92
$ws.Range("H92").Value = 2822.9412
$ws.Range("I92").Value = 1750
$ws.Range("J92").Value = 3153.077
$ws.Range("K92").Value = 5250
$ws.Range("L92").Value = 9459.231
$ws.Range("M92").Value = -4002
$ws.Range("N92").Value = -11955.231
# Row 94
$ws.Range("H94").Value = 5530.6
$ws.Range("I94").Value = 926
$ws.Range("J94").Value = 6681.75
$ws.Range("K94").Value = 2778
$ws.Range("L94").Value = 20045.25
$ws.Range("M94").Value = -2102
$ws.Range("N94").Value = -21397.25
# Row 98
$ws.Range("H98").Value = 2259.25
$ws.Range("I98").Value = 1848
$ws.Range("J98").Value = 2464.875
$ws.Range("K98").Value = 5544
$ws.Range("L98").Value = 7394.625
$ws.Range("M98").Value = -4046
$ws.Range("N98").Value = -10390.625
# Row 112
$ws.Range("H112").Value = 100005890
$ws.Range("I112").Value = 125006376
$ws.Range("J112").Value = 3925
$ws.Range("K112").Value = 375019128
$ws.Range("L112").Value = 11775
$ws.Range("M112").Value = -375018020
$ws.Range("N112").Value = -13991
# Row 128
$ws.Range("H128").Value = 134619
$ws.Range("I128").Value = 134619
$ws.Range("K128").Value = 403857
$ws.Range("M128").Value = -398877
# Row 134
$ws.Range("H134").Value = 6869.5
$ws.Range("I134").Value = 4151.75
$ws.Range("K134").Value = 12455.25
$ws.Range("M134").Value = -7385.25
# Row 139
$ws.Range("H139").Value = 3547.348
$ws.Range("I139").Value = 2119.7693
$ws.Range("K139").Value = 6359.3079
$ws.Range("M139").Value = -1219.3079
# Row 140
$ws.Range("H140").Value = 1530.174
$ws.Range("I140").Value = 1300
$ws.Range("K140").Value = 3900
$ws.Range("M140").Value = 1280

# ---------------------------------------------------------------
# GSM tab
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")
# Row 15
$ws.Range("H15").Value = 46247.625
$ws.Range("J15").Value = 46247.625
$ws.Range("L15").Value = 46247.625
$ws.Range("N15").Value = -46823.625
# Row 19
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()
# Row 21
$ws.Range("H21").Value = 22998.666
$ws.Range("I21").Value = 29498
$ws.Range("J21").Value = 10000
$ws.Range("K21").Value = 29498
$ws.Range("L21").Value = 10000
$ws.Range("M21").Value = -29325
$ws.Range("N21").Value = -10346
# Row 30
$ws.Range("H30").Value = 22998.666
$ws.Range("I30").Value = 29498
$ws.Range("J30").Value = 10000
$ws.Range("K30").Value = 29498
$ws.Range("L30").Value = 10000
$ws.Range("M30").Value = -29393
$ws.Range("N30").Value = -10210
# Row 81
$ws.Range("H81").Value = 46247.625
$ws.Range("J81").Value = 46247.625
$ws.Range("L81").Value = 46247.625
$ws.Range("N81").Value = -48243.625
# Row 84
$ws.Range("H84").Value = 46247.625
$ws.Range("J84").Value = 46247.625
$ws.Range("L84").Value = 138742.875
$ws.Range("N84").Value = -148726.875
# Row 102
$ws.Range("H102").Value = 2308.228
$ws.Range("I102").Value = 1627.5116
$ws.Range("K102").Value = 1627.5116
$ws.Range("M102").Value = -5.511600000000044
# Row 132
$ws.Range("H132").Value = 2099
$ws.Range("I132").Value = 1354.16
$ws.Range("J132").Value = 4426.625
$ws.Range("K132").Value = 4062.48
$ws.Range("L132").Value = 13279.875
$ws.Range("M132").Value = -1532.48
$ws.Range("N132").Value = -18339.875
# Row 136
$ws.Range("H136").Value = 60806.285
$ws.Range("J136").Value = 61440.668
$ws.Range("L136").Value = 184322.004
$ws.Range("N136").Value = -189422.004

# ---------------------------------------------------------------
# LTW tab
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")
# Row 13
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").ClearContents()
# Row 23
$ws.Range("H23").Value = 20006
$ws.Range("I23").Value = 20006
$ws.Range("K23").Value = 20006
$ws.Range("M23").Value = -19776
# Row 92
$ws.Range("H92").Value = 42142.855
$ws.Range("J92").Value = 42142.855
$ws.Range("L92").Value = 42142.855
$ws.Range("N92").Value = -47134.855

# ---------------------------------------------------------------
# WVR tab
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")
# Row 109
$ws.Range("H109").Value = 82500
$ws.Range("J109").Value = 82500
$ws.Range("L109").Value = 82500
$ws.Range("N109").Value = -85274
# Row 125
$ws.Range("H125").Value = 84087.875
$ws.Range("J125").Value = 84087.875
$ws.Range("L125").Value = 84087.875
$ws.Range("N125").Value = -93927.875
